$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# price-history rows down by one.
$ws.Rows("2:2").Insert()

# The newly inserted row 2 is blank; force the date-like text columns
# (A and E) to Text format first so Excel does not auto-convert the
# "dd-mm-yyyy" strings into date serial numbers when we assign them.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "14-11-2025"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 297.15
$ws.Range("E2").Value = "01-11-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# Re-apply the same cell formatting used by the rest of the data rows
# (copied from row 3, which holds what used to be row 2) so the new row
# matches the sheet's existing look exactly.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The row insert shifted every row's data down by one, but it left the
# worksheet's hyperlink anchors pointing at their old row numbers. Clear
# all of them out and recreate one hyperlink per data row (2-101),
# targeting each row's own Circular Link column text, so every F-cell's
# hyperlink matches the (now-shifted) value shown in it.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 101; $r++) {
  $url = $ws.Cells.Item($r, 6).Value2
  $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url) | Out-Null
}

$wb.Save()
